$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").NumberFormat = "@"
$meta.Range("B3").Value = "1.8.11"
$meta.Range("B8").Value = "2024-06-13T17:23:26-04:00"

# --- Elements sheet updates ---
$els = $wb.Worksheets.Item("Elements")

# Row 6 = Extension.value[x] ; Type(s) column K, Slicing Rules column AE
$els.Range("K6").Value = "base64Binary`nbooleancanonicalcodedatedateTimedecimalidinstantintegermarkdownoidpositiveIntstringtimeunsignedInturiurluuidAddressAgeAnnotationAttachmentCodeableConceptCodingContactPointCountDistanceDurationHumanNameIdentifierMoneyPeriodQuantityRangeRatioReferenceSampledDataSignatureTimingContactDetailContributorDataRequirementExpressionParameterDefinitionRelatedArtifactTriggerDefinitionUsageContextDosageMeta"
$els.Range("AE6").Value = "open"

# Column K width grows (content-driven autofit caps at 255)
$els.Columns.Item(11).ColumnWidth = 255
